$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Type the header row across A1:C1, as if the user typed each value
# and pressed Tab, leaving the selection on D1 afterwards.
$ws.Range("A1").Value = "month "
$ws.Range("B1").Value = "cgst"
$ws.Range("C1").Value = "sgst"

$ws.Range("D1").Select()
